# Renaming Target Variables, correcting Typo
#
# Sheet1 row 1 (AK1:AP1) holds the "target variable" headers for the
# juice-rating columns. Rename them and refresh the typo'd label.
#
#   AK1  Median            -> Med_Rating
#   AL1  Average           -> Avg_Rating
#   AM1  Rounded_Average   -> Avg_Rating_Rounded
#   AN1  Average_Recom     -> Recommondation
#   AO1  Avg>=6            -> Dummy_Avg
#   AP1  Med>=6            -> Dummy_Median

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AK1").Value = "Med_Rating"
$ws.Range("AL1").Value = "Avg_Rating"
$ws.Range("AM1").Value = "Avg_Rating_Rounded"
$ws.Range("AN1").Value = "Recommondation"
$ws.Range("AO1").Value = "Dummy_Avg"
$ws.Range("AP1").Value = "Dummy_Median"

# Restore the sheet's last active selection (cell AN10, near the renamed
# columns) to match the author's cursor position when they saved.
$ws.Range("AN10").Select()
